$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Name/Link swap for rows 20 and 21 (plain text, safe to set directly)
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"

# Force columns D:E to text format so numeric-looking strings are not
# auto-coerced into floating point numbers (preserves exact text + trailing zeros)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.373.86"
$ws.Range("D3").Value = "1.882.15"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "0.7118"
$ws.Range("D6").Value = "242.98"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "0.08025"
$ws.Range("E8").Value = "  +2.92%  "
$ws.Range("D9").Value = "0.3150"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "25.08"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").Value = "0.08335"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").Value = "1.909.11"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").Value = "5.268"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "95.03"
$ws.Range("E14").Value = "  +4.08%  "
$ws.Range("D15").Value = "0.7185"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "6.369"
$ws.Range("E16").Value = "  +5.20%  "
$ws.Range("D17").Value = "0.000008674"
$ws.Range("E17").Value = "  +5.38%  "
$ws.Range("D18").Value = "29.390.12"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "243.00"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "2.150.80"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").Value = "13.34"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "7.838"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "0.1574"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").Value = "9.096"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Value = "163.42"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "1.512"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "4.443"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "4.354"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "1.207"
$ws.Range("E32").Value = "  -6.22%  "
$ws.Range("D33").Value = "0.05397"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").Value = "1.948"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D35").Value = "0.7748"
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("D36").Value = "1.185"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "2.690"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").Value = "1.274.55"
$ws.Range("E39").Value = "  +3.94%  "
$ws.Range("D40").Value = "2.746"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").Value = "6.524"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").Value = "0.9196"
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("D43").Value = "113.09"
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("D44").Value = "74.49"
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("E46").Value = "  +4.56%  "
$ws.Range("D47").Value = "2.042.66"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("D48").Value = "1.815"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "0.5227"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "9.583"
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("D51").Value = "0.4377"
$ws.Range("E51").Value = "  +1.10%  "

# Restore original (default) cell style now that values are set as text
$ws.Range("D2:E51").Style = "Normal"
